# The sheet contains a long list of weekly price records for
# "Comercializadora del Agro de Limarí - Frutilla", grouped in blocks of
# three rows (Especial / Primera / Segunda) per reporting date, ordered
# from most-recent date (top) to oldest date (bottom).
#
# This edit adds one new, more-recent weekly record. Concretely, a new
# 3-row block is logically inserted right after the existing top block
# (row 378) and before the block that used to start at row 379 - pushing
# every following row down by one. The new block only has a single row
# of data supplied in the source diff (quality "Segunda"), reusing the
# other descriptive fields (market, region, product codes, unit, origin)
# from the block it displaces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at position 379; rows 379:435 shift down to 380:436
# and the sheet dimension grows from A1:T435 to A1:T436 automatically.
$ws.Rows("379:379").Insert()

# Seed the new row 379 with the same descriptive values as the row now
# sitting at 380 (the block that used to occupy row 379), then overwrite
# the fields that actually differ for the new record.
$src = $ws.Range("A380:T380")
$dst = $ws.Range("A379:T379")
$dst.Value2 = $src.Value2

$ws.Range("D379").Value2 = 44776
$ws.Range("L379").Value2 = "Segunda"
$ws.Range("M379").Value2 = 200
$ws.Range("N379").Value2 = 23000
$ws.Range("O379").Value2 = 24000
$ws.Range("P379").Value2 = 23500
$ws.Range("S379").Value2 = 3357
